$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.40"
$ws.Range("D2").ClearFormats()
$ws.Range("D4").Value = "'5.318"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").Value = "'0.05880"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = "'3.393"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'6.375"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").Value = "'0.8134"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = "'0.9573"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").Value = "'0.1418"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").Value = "'0.03728"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'0.07375"
$ws.Range("D12").ClearFormats()
$ws.Range("D14").Value = "'4.415"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = "'0.09400"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "'0.001592"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").Value = "'0.04811"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Value = "'0.0005905"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").Value = "'0.006100"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").Value = "'0.004080"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").Value = "'0.0009886"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").Value = "'0.00009710"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").Value = "'3.685"
$ws.Range("D23").ClearFormats()
$ws.Range("D26").Value = "'0.1275"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").Value = "'0.0002474"
$ws.Range("D27").ClearFormats()
$ws.Range("D40").Value = "'0.03897"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Value = "'0.006767"
$ws.Range("D41").ClearFormats()
$ws.Range("D43").Value = "'0.002703"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").Value = "'0.005919"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = "'0.00005677"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").ClearFormats()
$ws.Range("D47").Value = "'0.6521"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").Value = "'0.06580"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").Value = "'0.01011"
$ws.Range("D50").ClearFormats()
